$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest EUR->ARS quote as a new row (row 73).
# Column A holds a date-look-alike string ("2025-10-12"). Assigning it
# straight to .Value would make Excel auto-convert it into a date serial,
# so briefly force Text format, enter the value, then clear the format
# again so the cell ends up back on the sheet's default (General) style,
# exactly like the other rows, while keeping the literal text value.
$ws.Range("A73").NumberFormat = "@"
$ws.Range("A73").Value = "2025-10-12"
$ws.Range("A73").ClearFormats()

$ws.Range("B73").Value = "15:18:50"
$ws.Range("C73").Value = "1.00 EUR = 1,756.2048"
